# VirtualTradingLogs.xlsx - "Updated Logs for Test#5"
# Fills in the new trading-day column (I = "6 (April 3, 2023)") with the
# figures captured for that date, and refreshes the Totals / Gain% rows
# and the baseline legend to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- view state (best effort; cosmetic only) -----------------------------
$excel.ActiveWindow.Zoom = 60
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 4

# --- per-stock figures for the new date column (I) ------------------------
# Each source cell below already carries the correct "bought" (green) /
# "sold" (red) look for its row, so we copy its formatting across before
# writing the new number, exactly like the original author highlighting a
# matching neighbor cell and filling the new entry in the same style.

# I10 (TEL) - red/"sold" currency style, same look as E4
$ws.Range("E4").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 33125.85

# I14 (MER) - green/"bought" currency style, same look as C3
$ws.Range("C3").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 15525.67

# I16 (PGOLD) - red/"sold" currency style
$ws.Range("E4").Copy()
$ws.Range("I16").PasteSpecial(-4122)
$ws.Range("I16").Value = 16203.67

# I17 (LTG) - red/"sold" fill, but plain-number format this time
$ws.Range("D3").Copy()
$ws.Range("I17").PasteSpecial(-4122)
$ws.Range("I17").NumberFormat = "#,##0.00"
$ws.Range("I17").Value = 4897.16

# I18 (MPI) - red/"sold" currency style
$ws.Range("E4").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 18383.98

# I19 (AP) - red/"sold" currency style
$ws.Range("E4").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("I19").Value = 18284.87

# I20 (RRHI) - green/"bought" currency style
$ws.Range("C3").Copy()
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("I20").Value = 2740.31

$excel.CutCopyMode = $false

# --- Total Cash / Market Value for the new date ---------------------------
$ws.Range("I28").Value = 923293.89
$ws.Range("I29").Value = 68225.960000000006

# --- Total Equity row: re-enter the running totals so C30:H30 becomes a
#     shared fill (matches the rest of the row) and I30 gets its own sum ---
$ws.Range("C30:H30").Formula = "=SUM(C28:C29)"
$ws.Range("I30").Formula = "=SUM(I28:I29)"

# --- Baseline (PSEI) legend row: extend with the new date's figure --------
$ws.Range("H32").Copy()
$ws.Range("I32").PasteSpecial(-4122)
$ws.Range("I32").Value = -0.0072

$excel.CutCopyMode = $false

# --- restore the saved selection ------------------------------------------
$ws.Range("J35").Select()
